$wb = $excel.ActiveWorkbook
$nl = [char]10

foreach ($ws in $wb.Worksheets) {
    foreach ($suffix in @("[1]", "[2]", "[3]", "[4]", "[5]")) {
        $ws.Cells.Replace($suffix, "")
    }
    $ws.Cells.Replace($nl, " ")
}
